$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Journal entry for 2018-02-09 (serial 43140) - row 15
$ws.Range("A15").Value = 43140
$ws.Range("B15").Value = "absent du au décès de ma grand-maman"
$ws.Range("D15").Value = "1h30"

# Journal entry for 2018-02-12 (serial 43143) - row 16
$ws.Range("A16").Value = 43143
$ws.Range("B16").Value = "code classe Produit, Loueur, Location, et connection bd"
$ws.Range("D16").Value = "1h30"

# Update the active selection to D17, as left by the author after editing
$ws.Range("D17").Select()
